$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for Guayaba (Agrícola del Norte S.A.
# de Arica). Insert it as the new row 48, pushing the existing rows 48-70
# down to 49-71 (dimension grows from A1:T70 to A1:T71).
$ws.Rows.Item(48).Insert()

$ws.Cells.Item(48, 1).Value  = 1
$ws.Cells.Item(48, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(48, 4).Value  = 45097
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value  = 15
$ws.Cells.Item(48, 6).Value  = "Fruta"
$ws.Cells.Item(48, 7).Value  = 100108
$ws.Cells.Item(48, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(48, 9).Value  = 100108001
$ws.Cells.Item(48, 10).Value = "Guayaba"
$ws.Cells.Item(48, 11).Value = "Sin especificar"
$ws.Cells.Item(48, 12).Value = "Segunda"
$ws.Cells.Item(48, 13).Value = 200
$ws.Cells.Item(48, 14).Value = 5000
$ws.Cells.Item(48, 15).Value = 6000
$ws.Cells.Item(48, 16).Value = 5500
$ws.Cells.Item(48, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(48, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 19).Value = 550
$ws.Cells.Item(48, 20).Value = 10
